$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New pin map rows appended after the existing data (rows 2-19 already populated).
# Cells are written in the same order the original author entered them so the
# resulting shared-string table ordering matches exactly.
$ws.Cells.Item(20, 1).Value = "PROJECTILE_COMMS_ENABLE_OP"
$ws.Cells.Item(20, 2).Value = "PB21"

$ws.Cells.Item(21, 1).Value = "FILL_VALVE_COARSE_OP"
$ws.Cells.Item(21, 2).Value = "PB17"

$ws.Cells.Item(22, 2).Value = "PB16"
$ws.Cells.Item(23, 2).Value = "PD21"
$ws.Cells.Item(24, 2).Value = "PD20"

$ws.Cells.Item(22, 1).Value = "FILL_VALVE_FINE_OP"
$ws.Cells.Item(23, 1).Value = "DUMP_VALVE_COARSE_OP"
$ws.Cells.Item(24, 1).Value = "DUMP_VALVE_FINE_OP"

$ws.Cells.Item(25, 1).Value = "SYSTEM_CHARGE_ENABLE_FB"
$ws.Cells.Item(25, 2).Value = "PC23"

$ws.Cells.Item(26, 1).Value = "SYSTEM_PRESSURE_ENABLE_FB"
$ws.Cells.Item(26, 2).Value = "PC22"

$ws.Cells.Item(20, 3).Value = "GPIO out, default low."
$ws.Cells.Item(21, 3).Value = "GPIO out, default low."
$ws.Cells.Item(22, 3).Value = "GPIO out, default low."
$ws.Cells.Item(23, 3).Value = "GPIO out, default low."
$ws.Cells.Item(24, 3).Value = "GPIO out, default low."
$ws.Cells.Item(25, 3).Value = "GPIO in, pull up."
$ws.Cells.Item(26, 3).Value = "GPIO in, pull up."

$ws.Range("H27").Select()
